# Generate Report for Handback
# The localization run for the "68735751-..." source file failed during the
# handback transform (file-name mismatch), so the status row for that file
# switches from "Ready for handoff" to "Handback transform failed" on every
# sheet that surfaces it, and the per-locale "Error Detail" column gets the
# explanatory message. The Error Detail column is also widened so the new
# text is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 68735751-... file; columns E (zh-cn) and F (de-de)
# both held the old "Ready for handoff" status.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 (68735751-...) Status column (C) + Error Detail column (P)
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: zlfeihjs.wfn is different with handoff file name: 68735751-dfd8-46bf-ba8a-c8660daf5634.f381c23a69fb0b20a91177673f133ac38b2a92f0.zh-cn."

# de-de sheet: row 3 (68735751-...) Status column (C) + Error Detail column (P)
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: zlfeihjs.wfn is different with handoff file name: 68735751-dfd8-46bf-ba8a-c8660daf5634.f381c23a69fb0b20a91177673f133ac38b2a92f0.de-de."

# Widen the Error Detail column (P, the 16th column) on both locale sheets so
# the new error text displays. (39.17 "characters" round-trips to the OOXML
# col width of exactly 40, matching the other width="40" columns on this sheet.)
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
